# ===========================================================
# Edit script: restructure hyperparameter comparison workbook
#  - rename "Parameters" -> "Original dataset"
#  - add "Expanded dataset" and "Adam optimizer" sheets
#  - update titles, defined names, and selections
# ===========================================================

$wb = $excel.ActiveWorkbook

# --- Rename existing sheet and add two new sheets in order ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Original dataset"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Expanded dataset"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Adam optimizer"

# --- Update "Original dataset" title & header labels ---
$ws1.Range("B2").Value = "Comparison of hyperparameters for the feed-forward neural network (lr = 0.03)"
$ws1.Range("B4").Value = "#"
$ws1.Range("C4").Value = "epochs"

# --- "Expanded dataset" sheet: title + headers ---
$ws2.Range("B2").Value = "Comparison of hyperparameters for the feed-forward neural network (lr = 0.03, expanded dataset with rotated and flipped images)"
$ws2.Range("B2").Font.Bold = $true
$ws2.Range("B4").Value = "#"
$ws2.Range("C4").Value = "epochs"
$ws2.Range("D4").Value = "test weight_decay: 0.0, dropout: 0.0"
$ws2.Range("E4").Value = "train weight_decay: 0.0, dropout: 0.0"
$ws2.Range("F4").Value = "test weight_decay: 0.0, dropout: 0.3"
$ws2.Range("G4").Value = "train weight_decay: 0.0, dropout: 0.3"
$ws2.Range("H4").Value = "test weight_decay: 0.0, dropout: 0.5"
$ws2.Range("I4").Value = "train weight_decay: 0.0, dropout: 0.5"

# --- "Expanded dataset" sheet: data rows (B5:I14) ---
$ws2.Cells.Item(5, 2).Value = 0
$ws2.Cells.Item(5, 3).Value = 10
$ws2.Cells.Item(5, 4).Value = 0.82799999999999996
$ws2.Cells.Item(5, 5).Value = 0.86199999999999999
$ws2.Cells.Item(5, 6).Value = 0.81144444444444397
$ws2.Cells.Item(5, 7).Value = 0.83311111111111102
$ws2.Cells.Item(5, 8).Value = 0.80144444444444396
$ws2.Cells.Item(5, 9).Value = 0.81565079365079296
$ws2.Cells.Item(6, 2).Value = 1
$ws2.Cells.Item(6, 3).Value = 20
$ws2.Cells.Item(6, 4).Value = 0.84733333333333305
$ws2.Cells.Item(6, 5).Value = 0.93088888888888799
$ws2.Cells.Item(6, 6).Value = 0.84422222222222199
$ws2.Cells.Item(6, 7).Value = 0.898126984126984
$ws2.Cells.Item(6, 8).Value = 0.82699999999999996
$ws2.Cells.Item(6, 9).Value = 0.87431746031746005
$ws2.Cells.Item(7, 2).Value = 2
$ws2.Cells.Item(7, 3).Value = 30
$ws2.Cells.Item(7, 4).Value = 0.84355555555555495
$ws2.Cells.Item(7, 5).Value = 0.95499999999999996
$ws2.Cells.Item(7, 6).Value = 0.83622222222222198
$ws2.Cells.Item(7, 7).Value = 0.92280952380952297
$ws2.Cells.Item(7, 8).Value = 0.83077777777777695
$ws2.Cells.Item(7, 9).Value = 0.91007936507936504
$ws2.Cells.Item(8, 2).Value = 3
$ws2.Cells.Item(8, 3).Value = 40
$ws2.Cells.Item(8, 4).Value = 0.84211111111111103
$ws2.Cells.Item(8, 5).Value = 0.97466666666666602
$ws2.Cells.Item(8, 6).Value = 0.83855555555555505
$ws2.Cells.Item(8, 7).Value = 0.94658730158730098
$ws2.Cells.Item(8, 8).Value = 0.82811111111111102
$ws2.Cells.Item(8, 9).Value = 0.92312698412698402
$ws2.Cells.Item(9, 2).Value = 4
$ws2.Cells.Item(9, 3).Value = 50
$ws2.Cells.Item(9, 4).Value = 0.82511111111111102
$ws2.Cells.Item(9, 5).Value = 0.98295238095238002
$ws2.Cells.Item(9, 6).Value = 0.83288888888888801
$ws2.Cells.Item(9, 7).Value = 0.95223809523809499
$ws2.Cells.Item(9, 8).Value = 0.83
$ws2.Cells.Item(9, 9).Value = 0.94266666666666599
$ws2.Cells.Item(10, 2).Value = 5
$ws2.Cells.Item(10, 3).Value = 60
$ws2.Cells.Item(10, 4).Value = 0.82822222222222197
$ws2.Cells.Item(10, 5).Value = 0.98668253968253905
$ws2.Cells.Item(10, 6).Value = 0.82955555555555505
$ws2.Cells.Item(10, 7).Value = 0.96342857142857097
$ws2.Cells.Item(10, 8).Value = 0.83155555555555505
$ws2.Cells.Item(10, 9).Value = 0.95079365079364997
$ws2.Cells.Item(11, 2).Value = 6
$ws2.Cells.Item(11, 3).Value = 70
$ws2.Cells.Item(11, 4).Value = 0.82877777777777695
$ws2.Cells.Item(11, 5).Value = 0.99574603174603105
$ws2.Cells.Item(11, 6).Value = 0.83
$ws2.Cells.Item(11, 7).Value = 0.97231746031746003
$ws2.Cells.Item(11, 8).Value = 0.82211111111111101
$ws2.Cells.Item(11, 9).Value = 0.949031746031746
$ws2.Cells.Item(12, 2).Value = 7
$ws2.Cells.Item(12, 3).Value = 80
$ws2.Cells.Item(12, 4).Value = 0.822888888888888
$ws2.Cells.Item(12, 5).Value = 0.96053968253968203
$ws2.Cells.Item(12, 6).Value = 0.83088888888888801
$ws2.Cells.Item(12, 7).Value = 0.97649206349206297
$ws2.Cells.Item(12, 8).Value = 0.82199999999999995
$ws2.Cells.Item(12, 9).Value = 0.96077777777777695
$ws2.Cells.Item(13, 2).Value = 8
$ws2.Cells.Item(13, 3).Value = 90
$ws2.Cells.Item(13, 4).Value = 0.83277777777777695
$ws2.Cells.Item(13, 5).Value = 0.99995238095238004
$ws2.Cells.Item(13, 6).Value = 0.82177777777777705
$ws2.Cells.Item(13, 7).Value = 0.97319047619047605
$ws2.Cells.Item(13, 8).Value = 0.82877777777777695
$ws2.Cells.Item(13, 9).Value = 0.96706349206349196
$ws2.Cells.Item(14, 2).Value = 9
$ws2.Cells.Item(14, 3).Value = 100
$ws2.Cells.Item(14, 4).Value = 0.832666666666666
$ws2.Cells.Item(14, 5).Value = 0.99996825396825395
$ws2.Cells.Item(14, 6).Value = 0.82488888888888801
$ws2.Cells.Item(14, 7).Value = 0.98312698412698396
$ws2.Cells.Item(14, 8).Value = 0.82877777777777695
$ws2.Cells.Item(14, 9).Value = 0.96949206349206296

# --- "Expanded dataset" sheet: column widths ---
$ws2.Columns.Item(1).ColumnWidth = 1.0
$ws2.Columns.Item(2).ColumnWidth = 2.5
$ws2.Columns.Item(3).ColumnWidth = 5.833333333333333
$ws2.Columns.Item(4).ColumnWidth = 29.666666666666668
$ws2.Columns.Item(5).ColumnWidth = 30.5
$ws2.Columns.Item(6).ColumnWidth = 29.666666666666668
$ws2.Columns.Item(7).ColumnWidth = 30.5
$ws2.Columns.Item(8).ColumnWidth = 30.5
$ws2.Columns.Item(9).ColumnWidth = 31.0

# --- "Adam optimizer" sheet: title + headers ---
$ws3.Range("B2").Value = "Comparison of hyperparameters for the feed-forward neural network (lr = 0.03, expanded dataset with rotated and flipped images) with Adam optimizer and number of minibatches == 1000"
$ws3.Range("B2:C2").Font.Bold = $true
$ws3.Range("B4").Value = "#"
$ws3.Range("C4").Value = "epochs"
$ws3.Range("D4").Value = "test weight_decay: 0.0, dropout: 0.0"
$ws3.Range("E4").Value = "train weight_decay: 0.0, dropout: 0.0"
$ws3.Range("F4").Value = "test weight_decay: 0.0, dropout: 0.3"
$ws3.Range("G4").Value = "train weight_decay: 0.0, dropout: 0.3"
$ws3.Range("H4").Value = "test weight_decay: 0.0, dropout: 0.5"
$ws3.Range("I4").Value = "train weight_decay: 0.0, dropout: 0.5"

# --- "Adam optimizer" sheet: data rows (B5:I14) ---
$ws3.Cells.Item(5, 2).Value = 0
$ws3.Cells.Item(5, 3).Value = 10
$ws3.Cells.Item(5, 4).Value = 0.83588888888888802
$ws3.Cells.Item(5, 5).Value = 0.94531746031746
$ws3.Cells.Item(5, 6).Value = 0.83155555555555505
$ws3.Cells.Item(5, 7).Value = 0.91253968253968198
$ws3.Cells.Item(5, 8).Value = 0.82655555555555504
$ws3.Cells.Item(5, 9).Value = 0.89107936507936503
$ws3.Cells.Item(6, 2).Value = 1
$ws3.Cells.Item(6, 3).Value = 20
$ws3.Cells.Item(6, 4).Value = 0.83411111111111103
$ws3.Cells.Item(6, 5).Value = 0.96041269841269805
$ws3.Cells.Item(6, 6).Value = 0.83099999999999996
$ws3.Cells.Item(6, 7).Value = 0.93546031746031699
$ws3.Cells.Item(6, 8).Value = 0.82444444444444398
$ws3.Cells.Item(6, 9).Value = 0.91857142857142804
$ws3.Cells.Item(7, 2).Value = 2
$ws3.Cells.Item(7, 3).Value = 30
$ws3.Cells.Item(7, 4).Value = 0.83277777777777695
$ws3.Cells.Item(7, 5).Value = 0.97846031746031703
$ws3.Cells.Item(7, 6).Value = 0.82811111111111102
$ws3.Cells.Item(7, 7).Value = 0.94363492063492005
$ws3.Cells.Item(7, 8).Value = 0.82477777777777705
$ws3.Cells.Item(7, 9).Value = 0.93326984126984103
$ws3.Cells.Item(8, 2).Value = 3
$ws3.Cells.Item(8, 3).Value = 40
$ws3.Cells.Item(8, 4).Value = 0.83244444444444399
$ws3.Cells.Item(8, 5).Value = 0.97396825396825304
$ws3.Cells.Item(8, 6).Value = 0.83144444444444399
$ws3.Cells.Item(8, 7).Value = 0.95376190476190403
$ws3.Cells.Item(8, 8).Value = 0.82133333333333303
$ws3.Cells.Item(8, 9).Value = 0.94420634920634905
$ws3.Cells.Item(9, 2).Value = 4
$ws3.Cells.Item(9, 3).Value = 50
$ws3.Cells.Item(9, 4).Value = 0.83188888888888801
$ws3.Cells.Item(9, 5).Value = 0.98161904761904695
$ws3.Cells.Item(9, 6).Value = 0.82322222222222197
$ws3.Cells.Item(9, 7).Value = 0.95171428571428496
$ws3.Cells.Item(9, 8).Value = 0.82277777777777705
$ws3.Cells.Item(9, 9).Value = 0.94993650793650797
$ws3.Cells.Item(10, 2).Value = 5
$ws3.Cells.Item(10, 3).Value = 60
$ws3.Cells.Item(10, 4).Value = 0.82977777777777695
$ws3.Cells.Item(10, 5).Value = 0.97912698412698396
$ws3.Cells.Item(10, 6).Value = 0.823888888888888
$ws3.Cells.Item(10, 7).Value = 0.95931746031746001
$ws3.Cells.Item(10, 8).Value = 0.82155555555555504
$ws3.Cells.Item(10, 9).Value = 0.95242857142857096
$ws3.Cells.Item(11, 2).Value = 6
$ws3.Cells.Item(11, 3).Value = 70
$ws3.Cells.Item(11, 4).Value = 0.82933333333333303
$ws3.Cells.Item(11, 5).Value = 0.98411111111111105
$ws3.Cells.Item(11, 6).Value = 0.82555555555555504
$ws3.Cells.Item(11, 7).Value = 0.96560317460317402
$ws3.Cells.Item(11, 8).Value = 0.81899999999999995
$ws3.Cells.Item(11, 9).Value = 0.95587301587301499
$ws3.Cells.Item(12, 2).Value = 7
$ws3.Cells.Item(12, 3).Value = 80
$ws3.Cells.Item(12, 4).Value = 0.82411111111111102
$ws3.Cells.Item(12, 5).Value = 0.98577777777777698
$ws3.Cells.Item(12, 6).Value = 0.82566666666666599
$ws3.Cells.Item(12, 7).Value = 0.96820634920634896
$ws3.Cells.Item(12, 8).Value = 0.81299999999999994
$ws3.Cells.Item(12, 9).Value = 0.95630158730158699
$ws3.Cells.Item(13, 2).Value = 8
$ws3.Cells.Item(13, 3).Value = 90
$ws3.Cells.Item(13, 4).Value = 0.82899999999999996
$ws3.Cells.Item(13, 5).Value = 0.98561904761904695
$ws3.Cells.Item(13, 6).Value = 0.82322222222222197
$ws3.Cells.Item(13, 7).Value = 0.96276190476190404
$ws3.Cells.Item(13, 8).Value = 0.81722222222222196
$ws3.Cells.Item(13, 9).Value = 0.95950793650793598
$ws3.Cells.Item(14, 2).Value = 9
$ws3.Cells.Item(14, 3).Value = 100
$ws3.Cells.Item(14, 4).Value = 0.83322222222222198
$ws3.Cells.Item(14, 5).Value = 0.99012698412698397
$ws3.Cells.Item(14, 6).Value = 0.81655555555555503
$ws3.Cells.Item(14, 7).Value = 0.96915873015873
$ws3.Cells.Item(14, 8).Value = 0.81699999999999995
$ws3.Cells.Item(14, 9).Value = 0.96379365079364998

# --- "Adam optimizer" sheet: column widths ---
$ws3.Columns.Item(1).ColumnWidth = 2.0
$ws3.Columns.Item(2).ColumnWidth = 1.3333333333333333
$ws3.Columns.Item(3).ColumnWidth = 6.666666666666667
$ws3.Columns.Item(4).ColumnWidth = 30.166666666666668
$ws3.Columns.Item(5).ColumnWidth = 31.0
$ws3.Columns.Item(6).ColumnWidth = 30.166666666666668
$ws3.Columns.Item(7).ColumnWidth = 31.0
$ws3.Columns.Item(8).ColumnWidth = 30.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 31.0

# --- Defined name: comparison_1557179772.093524 (local to "Adam optimizer") ---
$rngAdam = $ws3.Range("B4:I14")
$null = $ws3.Names.Add("comparison_1557179772.093524", $rngAdam)

# --- Selections per sheet (applied before final activation) ---
$null = $ws1.Range("B1").Select()
$null = $ws2.Range("C4:C14").Select()

# --- Activate "Adam optimizer" as the active tab, with B5 selected ---
$ws3.Activate()
$null = $ws3.Range("B5").Select()

